$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table originally only tracked a single "upgrade 1" battery/PV
# configuration for each metric (Units, Total Nominal Capacity, Investment,
# Yearly O&M Cost). This adds "upgrade 2" and "upgrade 3" variants of each
# metric (part of the Salvage Value / Battery Replacement feature), so the
# table grows from 4 metric-rows to 12 metric-rows (rows 6-17), all reset to
# a placeholder value of 0 pending recalculation.

# Extend the bold/bordered/centered label formatting (already used by A2:A9)
# down to the newly added label rows A10:A17.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$labels = @(
  "Units at upgrade 1",
  "Units at upgrade 2",
  "Units at upgrade 3",
  "Total Nominal Capacity at upgrade 1",
  "Total Nominal Capacity at upgrade 2",
  "Total Nominal Capacity at upgrade 3",
  "Investment at upgrade 1",
  "Investment at upgrade 2",
  "Investment at upgrade 3",
  "Yearly O&M Cost at upgrade 1",
  "Yearly O&M Cost at upgrade 2",
  "Yearly O&M Cost at upgrade 3"
)

$row = 6
foreach ($label in $labels) {
    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $row = $row + 1
}

$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
